$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old last data row (row 6, "ECs -> Resolving-Mac") entirely —
# the other data rows (2-5) shift up conceptually one target-cluster slot
# and get recomputed TPM-based numbers below.
$ws.Rows(6).Delete()

# Row 2: ECs -> FAPs
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7349876666666667
$ws.Range("H2").Value = 2.204963
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.08219466666666668
$ws.Range("N2").Value = 0.246584
$ws.Range("O2").Value = 0.01551908118057144
$ws.Range("P2").Value = 0.01610836182747308
$ws.Range("Q2").Value = 0.06041206626577779
$ws.Range("R2").Value = 0.5437085963920001
$ws.Range("S2").Value = 0.01551908118057144
$ws.Range("T2").Value = 0.01610836182747308

# Row 3: FAPs -> Inflammatory-Mac
$ws.Range("D3").Value = "Inflammatory-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7349876666666667
$ws.Range("H3").Value = 2.204963
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.682868333333333
$ws.Range("N3").Value = 5.048605
$ws.Range("O3").Value = 0.3177404488678863
$ws.Range("P3").Value = 0.3298054864224351
$ws.Range("Q3").Value = 1.236887469623889
$ws.Range("R3").Value = 11.131987226615
$ws.Range("S3").Value = 0.3177404488678863
$ws.Range("T3").Value = 0.3298054864224351

# Row 4: Inflammatory-Mac -> MuSCs
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7349876666666667
$ws.Range("H4").Value = 2.204963
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.581259
$ws.Range("N4").Value = 1.162518
$ws.Range("O4").Value = 0.1097468482294607
$ws.Range("P4").Value = 0.0759427236761118
$ws.Range("Q4").Value = 0.427218196139
$ws.Range("R4").Value = 2.563309176834
$ws.Range("S4").Value = 0.1097468482294607
$ws.Range("T4").Value = 0.0759427236761118

# Row 5: MuSCs -> Resolving-Mac
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7349876666666667
$ws.Range("H5").Value = 2.204963
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.950039666666667
$ws.Range("N5").Value = 8.850119000000001
$ws.Range("O5").Value = 0.5569936217220816
$ws.Range("P5").Value = 0.5781434280739799
$ws.Range("Q5").Value = 2.168242771177445
$ws.Range("R5").Value = 19.514184940597
$ws.Range("S5").Value = 0.5569936217220816
$ws.Range("T5").Value = 0.5781434280739799
